$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3:E3").ClearContents()
$ws.Range("B4").Value = 0.4002772348174767
$ws.Range("C4").Value = 0.5431348647183788
$ws.Range("D4").Value = 0.7886513243237703
$ws.Range("E4").Value = 0.1983021219847731
$ws.Range("B5").Value = 0.2946906611286113
$ws.Range("C5").Value = 0.7509706769625255
$ws.Range("D5").Value = 0.6277136061112309
$ws.Range("E5").Value = 0.6820443385204394
$ws.Range("B6").Value = 0.2221659431358408
$ws.Range("C6").Value = 0.5873160119473396
$ws.Range("D6").Value = 0.2318650345370438
$ws.Range("E6").Value = 0.2233401124591544
$ws.Range("B7").Value = 0.05231955695830004
$ws.Range("C7").Value = 0.9163739752124588
$ws.Range("D7").Value = 0.460125088921839
$ws.Range("E7").Value = 0.3943265412458657
$ws.Range("B8").Value = 0.0119699462578432
$ws.Range("C8").Value = 0.9681531945047134
$ws.Range("D8").Value = 0.3850480984353858
$ws.Range("E8").Value = 0.7769529032977942
$ws.Range("B9").Value = 0.09596340882460686
$ws.Range("C9").Value = 0.8594424714340876
$ws.Range("D9").Value = 0.09025011991474496
$ws.Range("E9").Value = 0.1465547978935255
$ws.Range("B10").Value = 0.05480131625056961
$ws.Range("C10").Value = 0.7107469028220068
$ws.Range("D10").Value = 0.7390505873732331
$ws.Range("E10").Value = 0.2859386108804015
$ws.Range("B11").Value = 0.06422856648110331
$ws.Range("C11").Value = 0.7679801630257539
$ws.Range("D11").Value = 0.6079700923505065
$ws.Range("E11").Value = 0.325544513299205
$ws.Range("B12").Value = 0.07365581671163701
$ws.Range("C12").Value = 0.8252134232295009
$ws.Range("D12").Value = 0.4768895973277797
$ws.Range("E12").Value = 0.3651504157180084
$ws.Range("B13").Value = 0.08308306694217071
$ws.Range("C13").Value = 0.8824466834332478
$ws.Range("D13").Value = 0.345809102305053
$ws.Range("E13").Value = 0.4047563181368118
$ws.Range("B14").Value = 0.09251031717270441
$ws.Range("C14").Value = 0.9396799436369949
$ws.Range("D14").Value = 0.2147286072823262
$ws.Range("E14").Value = 0.4443622205556153
$ws.Range("B15").Value = 0.135165339552934
$ws.Range("C15").Value = 0.5607200016458863
$ws.Range("D15").Value = 0.4009313365275903
$ws.Range("E15").Value = 0.5796378777642874
$ws.Range("B16").Value = 0.6907932158480315
$ws.Range("C16").Value = 0.7573630362045051
$ws.Range("D16").Value = 0.7557610640945464
$ws.Range("E16").Value = 0.197494817332417
$ws.Range("B17").Value = 0.5897160456816568
$ws.Range("C17").Value = 0.8325245464970247
$ws.Range("D17").Value = 0.5422165271291751
$ws.Range("E17").Value = 0.3455607618129297
$ws.Range("B18").Value = 0.488638875515282
$ws.Range("C18").Value = 0.9076860567895444
$ws.Range("D18").Value = 0.3286719901638039
$ws.Range("E18").Value = 0.4936267062934425
$ws.Range("B19").Value = 0.3875617053489072
$ws.Range("C19").Value = 0.982847567082064
$ws.Range("D19").Value = 0.1151274531984326
$ws.Range("E19").Value = 0.6416926507739552
$ws.Range("B20").Value = 0.8862605938668301
$ws.Range("C20").Value = 0.7270351856955494
$ws.Range("D20").Value = 0.6499926272132465
$ws.Range("E20").Value = 0.3452828256819239
$ws.Range("B21").Value = 0.9475179683522807
$ws.Range("C21").Value = 0.6574841711212315
$ws.Range("D21").Value = 0.1477722063434623
$ws.Range("E21").Value = 0.2268406064517345
